# Refresh the cryptocurrency price/volume snapshot (scheduled GitHub Actions update).
# Price (column D) cells are stored as plain text in the source data (e.g. thousand-dot
# separated values like '3.187.13'), so pin NumberFormat to Text before writing each one —
# otherwise Excel's automatic type detection would silently convert plain decimal-looking
# values ("210.28", "32.00", ...) into numbers and normalize/strip their text formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "80.389.38"
$ws.Cells.Item(2, 5).Value = "  +4.84%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.189.68"
$ws.Cells.Item(3, 5).Value = "  +2.02%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.15%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "210.28"
$ws.Cells.Item(5, 5).Value = "  +5.10%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "628.33"
$ws.Cells.Item(6, 5).Value = "  +0.74%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.278"
$ws.Cells.Item(7, 5).Value = "  +27.55%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.999"
$ws.Cells.Item(8, 5).Value = "  -0.05%  "
$ws.Cells.Item(9, 5).Value = "  +5.44%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "3.187.46"
$ws.Cells.Item(10, 5).Value = "  +1.99%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.589"
$ws.Cells.Item(11, 5).Value = "  +20.71%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0000261"
$ws.Cells.Item(12, 5).Value = "  +28.35%  "
$ws.Cells.Item(13, 5).Value = "  +1.16%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "3.774.73"
$ws.Cells.Item(14, 5).Value = "  +1.72%  "
$ws.Cells.Item(15, 5).Value = "  +0.53%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "32.00"
$ws.Cells.Item(16, 5).Value = "  +5.57%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "80.386.46"
$ws.Cells.Item(17, 5).Value = "  +4.93%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "3.191.17"
$ws.Cells.Item(18, 5).Value = "  +1.71%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "14.29"
$ws.Cells.Item(19, 5).Value = "  +3.88%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "3.02"
$ws.Cells.Item(20, 5).Value = "  +9.42%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "442.54"
$ws.Cells.Item(21, 5).Value = "  +9.28%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "9.20"
$ws.Cells.Item(22, 5).Value = "  -0.36%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.24"
$ws.Cells.Item(23, 5).Value = "  +13.32%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "6.94"
$ws.Cells.Item(24, 5).Value = "  +6.12%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "3.348.25"
$ws.Cells.Item(25, 5).Value = "  +1.56%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "76.40"
$ws.Cells.Item(26, 5).Value = "  +3.19%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "4.71"
$ws.Cells.Item(27, 5).Value = "  +1.54%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "10.91"
$ws.Cells.Item(28, 5).Value = "  +4.00%  "
$ws.Cells.Item(29, 5).Value = "  +0.07%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0000123"
$ws.Cells.Item(30, 5).Value = "  +9.88%  "
$ws.Cells.Item(31, 5).Value = "  +0.20%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "8.98"
$ws.Cells.Item(32, 5).Value = "  +5.08%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "560.38"
$ws.Cells.Item(33, 5).Value = "  +8.57%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.45"
$ws.Cells.Item(34, 5).Value = "  -1.34%  "
$ws.Cells.Item(35, 5).Value = "  +13.20%  "
$ws.Cells.Item(36, 5).Value = "  +1.93%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "23.05"
$ws.Cells.Item(37, 5).Value = "  +6.51%  "
$ws.Cells.Item(38, 5).Value = "  +20.10%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.999"
$ws.Cells.Item(39, 5).Value = "  -0.15%  "
$ws.Cells.Item(40, 5).Value = "  +5.64%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "20.79"
$ws.Cells.Item(41, 5).Value = "  +3.57%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "162.84"
$ws.Cells.Item(42, 5).Value = "  -0.57%  "
$ws.Cells.Item(43, 5).Value = "  +5.96%  "
$ws.Cells.Item(44, 5).Value = "  -0.03%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "190.33"
$ws.Cells.Item(45, 5).Value = "  -2.69%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.73"
$ws.Cells.Item(46, 5).Value = "  +10.31%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.82"
$ws.Cells.Item(47, 5).Value = "  +5.26%  "
$ws.Cells.Item(48, 5).Value = "  -2.73%  "
$ws.Cells.Item(49, 5).Value = "  +0.43%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "42.84"
$ws.Cells.Item(50, 5).Value = "  +1.53%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "4.26"
$ws.Cells.Item(51, 5).Value = "  +5.19%  "
